$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.673.99"
Set-TextValue $ws.Range("E2") "  +1.48%  "

Set-TextValue $ws.Range("D3") "1.868.03"
Set-TextValue $ws.Range("E3") "  +0.56%  "

Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  -0.08%  "

Set-TextValue $ws.Range("D5") "235.87"
Set-TextValue $ws.Range("E5") "  +1.04%  "

Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.06%  "

Set-TextValue $ws.Range("D7") "0.4720"
Set-TextValue $ws.Range("E7") "  -0.44%  "

Set-TextValue $ws.Range("D8") "0.2764"
Set-TextValue $ws.Range("E8") "  +1.03%  "

Set-TextValue $ws.Range("D9") "0.06379"
Set-TextValue $ws.Range("E9") "  -0.67%  "

Set-TextValue $ws.Range("D10") "17.95"
Set-TextValue $ws.Range("E10") "  +10.48%  "

Set-TextValue $ws.Range("D11") "1.887.41"
Set-TextValue $ws.Range("E11") "  +1.47%  "

Set-TextValue $ws.Range("D12") "0.07456"
Set-TextValue $ws.Range("E12") "  +0.01%  "

Set-TextValue $ws.Range("D13") "4.976"
Set-TextValue $ws.Range("E13") "  -0.62%  "

Set-TextValue $ws.Range("D14") "85.21"
Set-TextValue $ws.Range("E14") "  -0.09%  "

Set-TextValue $ws.Range("D15") "0.6375"
Set-TextValue $ws.Range("E15") "  +0.85%  "

Set-TextValue $ws.Range("D16") "30.625.16"
Set-TextValue $ws.Range("E16") "  +1.48%  "

Set-TextValue $ws.Range("D17") "245.32"
Set-TextValue $ws.Range("E17") "  +6.15%  "

Set-TextValue $ws.Range("D18") "1.001"

Set-TextValue $ws.Range("D19") "12.84"
Set-TextValue $ws.Range("E19") "  +0.36%  "

Set-TextValue $ws.Range("D20") "0.000007409"
Set-TextValue $ws.Range("E20") "  +0.94%  "

Set-TextValue $ws.Range("D21") "1.001"
Set-TextValue $ws.Range("E21") "  -0.08%  "

Set-TextValue $ws.Range("D22") "4.967"
Set-TextValue $ws.Range("E22") "  -2.34%  "

Set-TextValue $ws.Range("D23") "6.086"
Set-TextValue $ws.Range("E23") "  +1.44%  "

Set-TextValue $ws.Range("D24") "9.390"
Set-TextValue $ws.Range("E24") "  +1.63%  "

Set-TextValue $ws.Range("D25") "164.31"
Set-TextValue $ws.Range("E25") "  -1.76%  "

Set-TextValue $ws.Range("D26") "18.34"
Set-TextValue $ws.Range("E26") "  +2.96%  "

Set-TextValue $ws.Range("D27") "1.898"
Set-TextValue $ws.Range("E27") "  +1.10%  "

Set-TextValue $ws.Range("D28") "0.1020"
Set-TextValue $ws.Range("E28") "  +2.07%  "

Set-TextValue $ws.Range("D29") "1.381"
Set-TextValue $ws.Range("E29") "  -0.31%  "

Set-TextValue $ws.Range("D30") "4.080"
Set-TextValue $ws.Range("E30") "  -1.78%  "

Set-TextValue $ws.Range("D31") "3.866"
Set-TextValue $ws.Range("E31") "  -1.44%  "

Set-TextValue $ws.Range("D32") "0.04937"
Set-TextValue $ws.Range("E32") "  +0.95%  "

Set-TextValue $ws.Range("D33") "1.157"
Set-TextValue $ws.Range("E33") "  +1.41%  "

Set-TextValue $ws.Range("D34") "0.7125"
Set-TextValue $ws.Range("E34") "  -0.85%  "

Set-TextValue $ws.Range("D35") "2.708"
Set-TextValue $ws.Range("E35") "  +0.27%  "

Set-TextValue $ws.Range("D36") "0.01919"
Set-TextValue $ws.Range("E36") "  +0.44%  "

Set-TextValue $ws.Range("D37") "2.690"
Set-TextValue $ws.Range("E37") "  +2.06%  "

Set-TextValue $ws.Range("D38") "0.8841"
Set-TextValue $ws.Range("E38") "  -2.00%  "

Set-TextValue $ws.Range("D39") "1.996"
Set-TextValue $ws.Range("E39") "  +1.48%  "

Set-TextValue $ws.Range("D40") "105.74"
Set-TextValue $ws.Range("E40") "  -0.22%  "

Set-TextValue $ws.Range("D41") "1.000"
Set-TextValue $ws.Range("E41") "  +0.02%  "

Set-TextValue $ws.Range("D42") "0.4112"
Set-TextValue $ws.Range("E42") "  +0.46%  "

Set-TextValue $ws.Range("D43") "5.564"
Set-TextValue $ws.Range("E43") "  +0.02%  "

Set-TextValue $ws.Range("D44") "7.342"
Set-TextValue $ws.Range("E44") "  +4.05%  "

Set-TextValue $ws.Range("D45") "62.35"
Set-TextValue $ws.Range("E45") "  +2.10%  "

Set-TextValue $ws.Range("D46") "0.1228"
Set-TextValue $ws.Range("E46") "  +2.33%  "

Set-TextValue $ws.Range("D47") "8.681"
Set-TextValue $ws.Range("E47") "  -0.14%  "

Set-TextValue $ws.Range("D48") "33.76"
Set-TextValue $ws.Range("E48") "  +1.99%  "

Set-TextValue $ws.Range("D49") "1.383"
Set-TextValue $ws.Range("E49") "  -1.09%  "

Set-TextValue $ws.Range("D50") "0.05564"
Set-TextValue $ws.Range("E50") "  -0.25%  "

Set-TextValue $ws.Range("D51") "0.3705"
Set-TextValue $ws.Range("E51") "  +0.22%  "
